$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: add duration value in column G
$ws.Range("G19").Value = 1.5

# Row 21: new coursework entry
$ws.Range("A21").Value = "EE270"
$ws.Range("B21").Value = "E"
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = "E"
$ws.Range("E21").Value = 60
$ws.Range("F21").Value = "Exam"

# Update selection to match the authored state
$ws.Range("G19").Select()
